$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 559.1667
$ws.Range("I28").Value = 378.33334
$ws.Range("J28").Value = 1101.6666
$ws.Range("K28").Value = 378.33334
$ws.Range("L28").Value = 1101.6666
$ws.Range("M28").Value = 106.66666
$ws.Range("N28").Value = -2071.6666

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 17242474
$ws.Range("I135").Value = 1007.0833
$ws.Range("J135").Value = 100001520
$ws.Range("K135").Value = 9063.7497
$ws.Range("L135").Value = 900013680
$ws.Range("M135").Value = -6528.7497
$ws.Range("N135").Value = -900018750

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2112.5647
$ws.Range("I138").Value = 1253.3513
$ws.Range("J138").Value = 2774.875
$ws.Range("K138").Value = 3760.0539
$ws.Range("L138").Value = 8324.625
$ws.Range("M138").Value = 1379.9461
$ws.Range("N138").Value = -18604.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1504.2667
$ws.Range("I2").Value = 1744.909
$ws.Range("J2").Value = 842.5
$ws.Range("K2").Value = 1744.909
$ws.Range("L2").Value = 842.5
$ws.Range("M2").Value = -1631.909
$ws.Range("N2").Value = -1068.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H109").Value = 40108.75
$ws.Range("J109").Value = 40108.75
$ws.Range("L109").Value = 40108.75
$ws.Range("N109").Value = -42882.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H112").Value = 37438.5
$ws.Range("J112").Value = 37438.5
$ws.Range("L112").Value = 37438.5
$ws.Range("N112").Value = -40392.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1504.2667
$ws.Range("I116").Value = 1744.909
$ws.Range("J116").Value = 842.5
$ws.Range("K116").Value = 1744.909
$ws.Range("L116").Value = 842.5
$ws.Range("M116").Value = 549.0909999999999
$ws.Range("N116").Value = -5430.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1504.2667
$ws.Range("I3").Value = 1744.909
$ws.Range("J3").Value = 842.5
$ws.Range("K3").Value = 1744.909
$ws.Range("L3").Value = 842.5
$ws.Range("M3").Value = -1630.909
$ws.Range("N3").Value = -1070.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H110").Value = 45377.8
$ws.Range("J110").Value = 45377.8
$ws.Range("L110").Value = 45377.8
$ws.Range("N110").Value = -53557.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H119").Value = 45296
$ws.Range("J119").Value = 45296
$ws.Range("L119").Value = 45296
$ws.Range("N119").Value = -54972

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H120").Value = 46757
$ws.Range("J120").Value = 46757
$ws.Range("L120").Value = 46757
$ws.Range("N120").Value = -56433

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H122").Value = 40719.4
$ws.Range("J122").Value = 40719.4
$ws.Range("L122").Value = 40719.4
$ws.Range("N122").Value = -50519.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3196.85
$ws.Range("I134").Value = 1916.4166
$ws.Range("J134").Value = 4050.4722
$ws.Range("K134").Value = 5749.2498
$ws.Range("L134").Value = 12151.4166
$ws.Range("M134").Value = -3214.2498
$ws.Range("N134").Value = -17221.4166

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 421.33334
$ws.Range("I107").Value = 397.13043
$ws.Range("J107").Value = 464.15384
$ws.Range("K107").Value = 397.13043
$ws.Range("L107").Value = 464.15384
$ws.Range("M107").Value = 1522.86957
$ws.Range("N107").Value = -4304.15384

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H116").Value = 48491.668
$ws.Range("J116").Value = 48491.668
$ws.Range("L116").Value = 48491.668
$ws.Range("N116").Value = -57669.668

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 7577.5356
$ws.Range("I107").Value = 10418.6
$ws.Range("J107").Value = 5999.1665
$ws.Range("K107").Value = 31255.8
$ws.Range("L107").Value = 17997.4995
$ws.Range("M107").Value = -29335.8
$ws.Range("N107").Value = -21837.4995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1278.5555
$ws.Range("I122").Value = 1001.75
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 3005.25
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -555.25
$ws.Range("N122").Value = -9400

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2954.0435
$ws.Range("I7").Value = 2114
$ws.Range("J7").Value = 5334.1665
$ws.Range("K7").Value = 2114
$ws.Range("L7").Value = 5334.1665
$ws.Range("M7").Value = -2002
$ws.Range("N7").Value = -5558.1665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H110").Value = 38291.145
$ws.Range("J110").Value = 38291.145
$ws.Range("L110").Value = 38291.145
$ws.Range("N110").Value = -46471.145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H121").Value = 43412
$ws.Range("J121").Value = 43412
$ws.Range("L121").Value = 43412
$ws.Range("N121").Value = -46906

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2954.0435
$ws.Range("I126").Value = 2114
$ws.Range("J126").Value = 5334.1665
$ws.Range("K126").Value = 6342
$ws.Range("L126").Value = 16002.4995
$ws.Range("M126").Value = -3872
$ws.Range("N126").Value = -20942.4995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 18570
$ws.Range("I70").Value = 9000
$ws.Range("J70").Value = 20165
$ws.Range("K70").Value = 9000
$ws.Range("L70").Value = 20165
$ws.Range("M70").Value = -8685
$ws.Range("N70").Value = -20795

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H73").Value = 18570
$ws.Range("I73").Value = 9000
$ws.Range("J73").Value = 20165
$ws.Range("K73").Value = 9000
$ws.Range("L73").Value = 20165
$ws.Range("M73").Value = -7908
$ws.Range("N73").Value = -22349

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 6667450.5
$ws.Range("I107").Value = 863.6667
$ws.Range("J107").Value = 11111842
$ws.Range("K107").Value = 2591.0001
$ws.Range("L107").Value = 33335526
$ws.Range("M107").Value = -671.0001000000002
$ws.Range("N107").Value = -33339366

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H118").Value = 41995
$ws.Range("J118").Value = 41995
$ws.Range("L118").Value = 41995
$ws.Range("N118").Value = -45309

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 52857804
$ws.Range("I122").Value = 58730780
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 176192340
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -176189890
$ws.Range("N122").Value = -7900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2102039.8
$ws.Range("I126").Value = 2452279.8
$ws.Range("J126").Value = 600
$ws.Range("K126").Value = 7356839.399999999
$ws.Range("L126").Value = 1800
$ws.Range("M126").Value = -7354369.399999999
$ws.Range("N126").Value = -6740

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2171.6875
$ws.Range("I132").Value = 1766.5834
$ws.Range("J132").Value = 3387
$ws.Range("K132").Value = 5299.7502
$ws.Range("L132").Value = 10161
$ws.Range("M132").Value = -2769.7502
$ws.Range("N132").Value = -15221
